$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new diary entry row 17
$ws.Range("B17").Value = "Restrukturierung, Bereineinigung von Code-Redundanzen, Stage 1 von x"
$ws.Range("E17").Value = 1.5

# Copy style from the row above (E16) so the new value cell matches formatting
$ws.Range("E16").Copy()
$ws.Range("E17").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selected cell / range shown when the file is opened
$ws.Range("E18").Select()
